$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Speaking")

$ws.Range("A2").Value = "North Region "
$ws.Range("B2").NumberFormat = "m/d/yy"
$ws.Range("B2").Value = Get-Date -Year 2017 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "005774"
$ws.Range("D2").Value = 990801
